$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New values for rows 2-10, columns B-G
$data = @{
    2  = @(-0.0253360835587989, 1.135698984415373, 7.612000487694052, 2.758985409112207, 2.78578550780152, 52)
    3  = @(0.08600705238178873, 1.123115432295725, 5.609313565901017, 2.368398945680608, 2.390387976856151, 51)
    4  = @(0.0363203761572407, 0.9693363125750792, 4.390696811956634, 2.095398962478658, 2.116354596420031, 50)
    5  = @(0.1260652424010658, 1.03221972966461, 4.978658809889423, 2.231290839377382, 2.25081259337064, 49)
    6  = @(0.08428106382731648, 1.007929901751503, 5.095203452855081, 2.257255734925726, 2.27955206241149, 48)
    7  = @(0.09619512034022865, 1.131065357075751, 5.749733647272278, 2.397860222630226, 2.42991643159901, 36)
    8  = @(0.1381149580396931, 1.187418322474588, 5.903068426619335, 2.4296231038207, 2.461107737049912, 35)
    9  = @(0.1382699113425551, 1.562389918879535, 9.831370974918348, 3.135501710240061, 3.223266002762909, 18)
    10 = @(-0.7150897714758779, 1.341686922679743, 6.639370992174167, 2.576697691265735, 2.596308797503721, 11)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Range("B$row").Value = $values[0]
    $ws.Range("C$row").Value = $values[1]
    $ws.Range("D$row").Value = $values[2]
    $ws.Range("E$row").Value = $values[3]
    $ws.Range("F$row").Value = $values[4]
    $ws.Range("G$row").Value = $values[5]
}
